# Deploy: Update site content, preserve CNAME, add .nojekyll
# Update the Metadata sheet (Date + Count) and append new Concepts rows.

$wb = $excel.ActiveWorkbook
$metaWs = $wb.Worksheets.Item("Metadata")
$conceptsWs = $wb.Worksheets.Item("Concepts")

# ---- New concept rows (Code, Display, Definition); Level is always "1" ----
$newRows = @(
  @("what-matters", "What Matters to Me", "Patient-identified priorities and non-clinical goals"),
  @("patient-story", "Patient Story", "Narrative summary of patient background, preferences, and autobiography"),
  @("relational-engagement", "Relational Engagement Score", "Assessment of the quality of nurse-patient engagement (1-5)"),
  @("cfs-score", "Clinical Frailty Scale Score", "Total score for Rockwood Clinical Frailty Scale"),
  @("cfs-1", "Very Fit", "Robus, active, energetic and motivated"),
  @("cfs-2", "Well", "No active disease symptoms but less fit than category 1"),
  @("cfs-3", "Managing Well", "Medical problems are well controlled, but not regularly active"),
  @("cfs-4", "Vulnerable", "Not dependent for daily help, but symptoms limit activities"),
  @("cfs-5", "Mildly Frail", "Need help with high order IADLs (finances, transportation, heavy housework)"),
  @("cfs-6", "Moderately Frail", "Need help with all outside activities and some housekeeping"),
  @("cfs-7", "Severely Frail", "Completely dependent for personal care"),
  @("cfs-8", "Very Severely Frail", "Completely dependent, approaching end of life"),
  @("cfs-9", "Terminally Ill", "Approaching the end of life (life expectancy <6 months)"),
  @("4at-score", "4AT Delirium Score", "Total score for 4AT assessment"),
  @("4at-alertness", "Alertness", "4AT Item 1: Alertness"),
  @("4at-amt4", "AMT4 Score", "4AT Item 2: Abbreviated Mental Test 4"),
  @("4at-attention", "Attention", "4AT Item 3: Attention (Months Backwards)"),
  @("4at-acute-change", "Acute Change", "4AT Item 4: Acute Change or Fluctuating Course"),
  @("reasonable-adjustment", "Reasonable Adjustment", "Requirement for adjustment to care delivery (Equality Act)"),
  @("mca-assessment", "Mental Capacity Assessment", "Assessment of capacity to make a specific decision"),
  @("capacity-present", "Capacity Present", "Patient has capacity for this decision"),
  @("capacity-absent", "Capacity Absent", "Patient lacks capacity for this decision"),
  @("best-interest", "Best Interest Decision", "Decision made in patient's best interest"),
)

$startRow = 31 + 1
$formatSrc = $conceptsWs.Range("A31:D31")

$r = $startRow
foreach ($row in $newRows) {
  $dst = $conceptsWs.Range("A" + $r + ":D" + $r)
  $formatSrc.Copy()
  $dst.PasteSpecial(-4122)

  $conceptsWs.Cells.Item($r, 1).Value = "'1"
  $conceptsWs.Cells.Item($r, 2).Value = $row[0]
  $conceptsWs.Cells.Item($r, 3).Value = $row[1]
  $conceptsWs.Cells.Item($r, 4).Value = $row[2]

  # Re-apply the source formatting so the quote-prefix entry used to force
  # the Level column to stay text-typed doesn't leave behind a stray style.
  $formatSrc.Copy()
  $conceptsWs.Range("A" + $r).PasteSpecial(-4122)

  $r = $r + 1
}

# ---- Metadata updates ----
$metaWs.Cells.Item(8, 2).Value = "2025-12-26T15:22:58+00:00"

# Count ("53") looks numeric, so force it to stay text-typed (matching the
# original "30" string cell) via the quote-prefix trick, then restore the
# untouched formatting from a sibling cell that shares the same style.
$countFormatSrc = $metaWs.Range("B20")
$metaWs.Cells.Item(21, 2).Value = "'53"
$countFormatSrc.Copy()
$metaWs.Range("B21").PasteSpecial(-4122)
